$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.234.95"
$ws.Range("D3").Value = "1.845.45"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6649"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.87"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07475"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07759"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "1.850.25"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.196"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008770"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").Value = "29.200.23"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "2.097.85"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.222"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.9999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.39"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.649"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1405"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.515"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.142"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.057"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.192"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05375"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.854"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7476"
$ws.Range("D36").ClearFormats()
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.644"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").Value = "1.301.72"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01799"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.755"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.415"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +8.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9097"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.08259"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "1.996.21"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.757"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.87%  "

Write-Host "Applied crypto list update"
